$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# Paragraph 1: "Use Case:" <tab> -> drop the _GoBack bookmark, append "Delete Vehicle" run
$p1xml = '<w:body><w:p w:rsidR="00D5269F" w:rsidRDefault="00D5269F" w:rsidP="00D5269F"><w:r><w:t>Use Case:</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>Delete Vehicle</w:t></w:r></w:p></w:body>'
Replace-ParagraphXml 1 $p1xml

# Paragraph 2: "Precondition: " -> append sentence
$p2xml = '<w:body><w:p w:rsidR="00D5269F" w:rsidRDefault="00D5269F" w:rsidP="00F407D3"><w:r><w:t xml:space="preserve">Precondition: </w:t></w:r><w:r><w:t>User will have a list of current vehicles displayed on screen</w:t></w:r></w:p></w:body>'
Replace-ParagraphXml 2 $p2xml

# Paragraph 6: "Actors: " -> append "User"
$p6xml = '<w:body><w:p w:rsidR="00D5269F" w:rsidRDefault="00D5269F" w:rsidP="00F407D3"><w:r><w:t xml:space="preserve">Actors: </w:t></w:r><w:r><w:t>User</w:t></w:r></w:p></w:body>'
Replace-ParagraphXml 6 $p6xml

# Paragraph 10: "Description: " -> append full description incl. proofed "api"
$p10xml = '<w:body><w:p w:rsidR="00D5269F" w:rsidRDefault="00D5269F" w:rsidP="00F407D3"><w:r><w:t xml:space="preserve">Description: </w:t></w:r><w:r><w:t xml:space="preserve">User will have a list of the current vehicles on screen, they will then navigate to the part where the current vehicles are and then press on the delete button which will send a delete request to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and delete that vehicle from the list.</w:t></w:r></w:p></w:body>'
Replace-ParagraphXml 10 $p10xml

# Paragraph 30: "Post Condition:" -> append sentence, and move the _GoBack bookmark here (at the end)
$p30xml = '<w:body><w:p w:rsidR="00D5269F" w:rsidRDefault="00D5269F" w:rsidP="00F407D3"><w:r><w:t>Post Condition:</w:t></w:r><w:r><w:t xml:space="preserve"> User will no longer be able to see the vehicle in the list of vehicles from the api.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>'
Replace-ParagraphXml 30 $p30xml

Write-Output "done"
